$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fixture/prediction data for rows 2-19 (columns A-G)
# Columns: A=Fixture, B=Pick, C=Competition, D=Time, E=Win Info, F=Confidence %, G=Odds
$data = @(
    @("Liverpool vs Aston Villa", "Liverpool", "England", "1st Nov 20:00", "64 of 96 Tips", "67", "1.69"),
    @("Monaco vs Paris FC", "Monaco", "France", "1st Nov 18:00", "61 of 71 Tips", "86", "1.67"),
    @("Napoli vs Como", "Napoli", "Italy", "1st Nov 17:00", "56 of 74 Tips", "76", "1.95"),
    @("Real Madrid vs Valencia", "Real Madrid", "Spain", "1st Nov 20:00", "52 of 56 Tips", "93", "1.18"),
    @("Cremonese vs Juventus", "Juventus", "Italy", "1st Nov 19:45", "51 of 61 Tips", "84", "1.52"),
    @("Auxerre vs Marseille", "Marseille", "France", "1st Nov 20:05", "51 of 61 Tips", "84", "1.68"),
    @("Bayern Munich vs Bayer Leverkusen", "Bayern Munich", "Germany", "1st Nov 17:30", "29 of 34 Tips", "85", "1.23"),
    @("Almeria vs SD Eibar", "Almeria", "Spain", "1st Nov 20:00", "25 of 35 Tips", "71", "1.73"),
    @("Huachipato vs Universidad de Chile", "Universidad de Chile", "Chile", "2nd Nov 15:30", "25 of 26 Tips", "96", "1.73"),
    @("Real Madrid vs Valencia", "1st Half Result: Real Madrid", "Spain", "1st Nov 20:00", "23 of 29 Tips", "79", "1.52"),
    @("West Ham vs Newcastle", "Newcastle", "England", "2nd Nov 14:00", "23 of 27 Tips", "85", "1.61"),
    @("Cruzeiro vs Vitoria BA", "Cruzeiro", "Brazil", "1st Nov 19:00", "22 of 23 Tips", "96", "1.44"),
    @("Oriente Petrolero vs Guabira", "Oriente Petrolero", "Bolivia", "3rd Nov 19:00", "20 of 20 Tips", "100", "1.68"),
    @("Velez Sarsfield vs Talleres de Cordoba", "Velez Sarsfield", "Argentina", "1st Nov 20:00", "19 of 24 Tips", "79", "2.15"),
    @("Palestino vs Deportes Limache", "Palestino", "Chile", "3rd Nov 22:00", "19 of 19 Tips", "100", "1.79"),
    @("Bayern Munich vs Bayer Leverkusen", "1st Half Result: Bayern Munich", "Germany", "1st Nov 17:30", "17 of 23 Tips", "74", "1.57"),
    @("Club Brugge vs Dender", "Club Brugge", "Belgium", "1st Nov 17:15", "17 of 17 Tips", "100", "1.28"),
    @("Liverpool vs Aston Villa", "Double Chance: X2", "England", "1st Nov 20:00", "16 of 23 Tips", "70", "2.23")
)

$startRow = 2
$endRow = $startRow + $data.Length - 1

# Ensure the Confidence % and Odds columns stay text (as in the source file)
$ws.Range("F$startRow`:G$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
